# Hourly dataset refresh: updated regression coefficients (Coef./Std.Err./t/
# P>|t|/CI bounds/coef_pos) for the cap_gen_year12final sheet, including two
# rows (6 and 10) whose Std.Err./t/P/CI columns are no longer estimated
# (cleared) and two rows (18 and 25) that now have those columns populated.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.1065553440491354
$ws.Range("B3").Value = 0.1073220220072553
$ws.Range("H3").Value = 0.2138773660563906
$ws.Range("B4").Value = 0.09745546187922861
$ws.Range("H4").Value = 0.204010805928364
$ws.Range("B5").Value = 0.0749393821291105
$ws.Range("C5").Value = 0.00436756508448804
$ws.Range("D5").Value = 15.44556680224815
$ws.Range("E5").Value = 0.04111302272874021
$ws.Range("F5").Value = 0.06635261590659437
$ws.Range("G5").Value = 0.08352614835162808
$ws.Range("H5").Value = 0.1814947261782459
$ws.Range("B6").Value = 0.09806475263741116
$ws.Range("H6").Value = 0.2046200966865465
$ws.Range("B7").Value = 0.06849091206033495
$ws.Range("C7").Value = 0.002193589579185852
$ws.Range("D7").Value = 12.79972308625368
$ws.Range("E7").Value = 0.0101215292127131
$ws.Range("F7").Value = 0.06418917671874072
$ws.Range("G7").Value = 0.07279264740192995
$ws.Range("H7").Value = 0.1750462561094703
$ws.Range("B8").Value = 0.04936317317196443
$ws.Range("C8").Value = 0.00179921928313599
$ws.Range("D8").Value = 9.896076031876499
$ws.Range("E8").Value = 0.008413415385360275
$ws.Range("F8").Value = 0.04583505793116451
$ws.Range("G8").Value = 0.05289128841276468
$ws.Range("H8").Value = 0.1559185172210998
$ws.Range("B9").Value = 0.04906988299974192
$ws.Range("C9").Value = 0.002040482625462532
$ws.Range("D9").Value = 8.816382794271099
$ws.Range("E9").Value = 0.005404991362256598
$ws.Range("F9").Value = 0.04506104525694603
$ws.Range("G9").Value = 0.0530787207425381
$ws.Range("H9").Value = 0.1556252270488773
$ws.Range("B10").Value = 0.0499486100949007
$ws.Range("H10").Value = 0.1565039541440361
$ws.Range("B11").Value = 0.03272506821371153
$ws.Range("H11").Value = 0.1392804122628469
$ws.Range("B12").Value = 0.05117392714232886
$ws.Range("H12").Value = 0.1577292711914642
$ws.Range("B13").Value = 0.06807052312645226
$ws.Range("H13").Value = 0.1746258671755876
$ws.Range("B14").Value = 0.07726985812573917
$ws.Range("H14").Value = 0.1838252021748745
$ws.Range("B15").Value = 0.08069832767537394
$ws.Range("H15").Value = 0.1872536717245093
$ws.Range("B16").Value = 0.08364221506605111
$ws.Range("H16").Value = 0.1901975591151865
$ws.Range("B17").Value = 0.08505360333591129
$ws.Range("H17").Value = 0.1916089473850466
$ws.Range("B18").Value = -0.1065553440491354
$ws.Range("C18").Value = 0.008791455712720293
$ws.Range("D18").Value = -19.91194238406194
$ws.Range("E18").Value = 0.02634267328700736
$ws.Range("F18").Value = -0.1238444312855709
$ws.Range("G18").Value = -0.0892662568126998
$ws.Range("B19").Value = 0.08784019235903448
$ws.Range("H19").Value = 0.1943955364081698
$ws.Range("B20").Value = 0.08970924313437535
$ws.Range("H20").Value = 0.1962645871835107
$ws.Range("B21").Value = 0.0922956501447684
$ws.Range("H21").Value = 0.1988509941939038
$ws.Range("B22").Value = 0.09716712440069303
$ws.Range("H22").Value = 0.2037224684498284
$ws.Range("B23").Value = 0.1017495658785929
$ws.Range("H23").Value = 0.2083049099277283
$ws.Range("B24").Value = 0.1034528816909316
$ws.Range("H24").Value = 0.210008225740067
$ws.Range("B25").Value = 0.104008506325195
$ws.Range("C25").Value = 0.007342672136582881
$ws.Range("D25").Value = -464258035.6935544
$ws.Range("E25").Value = 0.04584449633748953
$ws.Range("F25").Value = 0.08957773653096798
$ws.Range("G25").Value = 0.118439276119421
$ws.Range("H25").Value = 0.2105638503743303
$ws.Range("B26").Value = 0.1085442960416581
$ws.Range("C26").Value = 0.00716812687958495
$ws.Range("D26").Value = 522695178.3343856
$ws.Range("E26").Value = 0.04100148843186305
$ws.Range("F26").Value = 0.0944616594714633
$ws.Range("G26").Value = 0.1226269326118539
$ws.Range("H26").Value = 0.2150996400907934
$ws.Range("B27").Value = 0.1121338848236663
$ws.Range("C27").Value = 0.007141991707207579
$ws.Range("D27").Value = 25.0630326955693
$ws.Range("E27").Value = 0.04587700108548545
$ws.Range("F27").Value = 0.09810172851232626
$ws.Range("G27").Value = 0.1261660411350068
$ws.Range("H27").Value = 0.2186892288728016
$ws.Range("B28").Value = 0.1139243743438063
$ws.Range("C28").Value = 0.00743520263279721
$ws.Range("D28").Value = 431441629.5995058
$ws.Range("E28").Value = 0.07278981944472064
$ws.Range("F28").Value = 0.09931268336114869
$ws.Range("G28").Value = 0.1285360653264641
$ws.Range("H28").Value = 0.2204797183929417
$ws.Range("B29").Value = 0.05238793070310608
$ws.Range("C29").Value = 0.001886243714204461
$ws.Range("D29").Value = 10.29448615874023
$ws.Range("E29").Value = 0.006560674141669946
$ws.Range("F29").Value = 0.04868089229358762
$ws.Range("G29").Value = 0.05609496911262121
$ws.Range("H29").Value = 0.1589432747522414

$ws.Range("C6:G6").ClearContents()
$ws.Range("C10:G10").ClearContents()
